$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

$ws.Range("A4").Value = "GitHub'ta Proje Oluşturma ve Push Etme"
$ws.Range("B4").Value = "Öznur URFAN"

$ws.Range("A4").Select()
